$d = $word.ActiveDocument

foreach ($para in $d.Paragraphs) {
    $rng = $para.Range
    if ($rng.Text -like "*Make .csv for CO2 conversions*") {
        # Strike through the paragraph mark (end of paragraph) and the text run
        $para.Range.Font.StrikeThrough = $true
        break
    }
}
